# Updated cryptos list on Sat Mar 16 07:13:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.258.15'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '3.729.08'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.17'
$ws.Range("E5").Value = '  +5.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.74'
$ws.Range("E6").Value = '  +7.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.726'
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.10'
$ws.Range("E11").Value = '  +11.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000292'
$ws.Range("E12").Value = '  -3.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.66'
$ws.Range("D14").Value = '4.314.94'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").Value = '3.720.22'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.44'
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.15'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.127'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("D20").Value = '69.065.96'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '413.29'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.59'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '89.72'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.07'
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.91'
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.93'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.82'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.05'
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  +0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.25'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.66'
$ws.Range("E31").Value = '  -4.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.83'
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("E33").Value = '  +3.70%  '
$ws.Range("E34").Value = '  +4.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '633.29'
$ws.Range("E35").Value = '  +3.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.70'
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").Value = '0.0₃0835'
$ws.Range("E37").Value = '  -9.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.415'
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("E42").Value = '  -1.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0448'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("E45").Value = '  +2.95%  '
$ws.Range("D46").Value = '2.883.88'
$ws.Range("E46").Value = '  +5.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.21'
$ws.Range("E47").Value = '  -2.25%  '
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.10'
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.51'
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  -19.93%  '
